$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Idrottare"
$ws.Range("B5").Value = "FYSISK"

[void]$ws.Range("B5").Select()
